# Fix the label order of caps
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: 20pF capacitors - relabel from C3,C6 to C3,C4
$ws.Range("E6").Value = "C3, C4"

# Row 5: 0.1uF capacitors - relabel from C1,C2,C7 to C1,C2,C5
$ws.Range("E5").Value = "C1, C2, C5"

# Row 4: 1uF capacitors - relabel from C4,C5,C10 to C6,C7,C10
$ws.Range("E4").Value = "C6, C7, C10"

# Bump the version note to reflect the fix
$ws.Range("B26").Value = "version 3.0.2"
